$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Build the new "Commit 8" block (rows 139-154), mirroring the layout of
# --- the existing "Commit 7" block (rows 122-137) directly above it (offset
# --- of +17 rows). Formats are copied cell-by-cell (not as a rectangular
# --- range) from the matching source cell so that rows/cells that have no
# --- content in the source block stay completely absent in the new block
# --- too, instead of materializing as blank placeholder cells.

$cellMap = @(
  @("A122","A139"),
  @("A123","A140"), @("B123","B140"), @("C123","C140"), @("D123","D140"),
  @("A125","A142"),
  @("A126","A143"), @("B126","B143"), @("D126","D143"), @("E126","E143"), @("F126","F143"),
  @("A127","A144"), @("B127","B144"), @("D127","D144"), @("E127","E144"), @("F127","F144"),
  @("A128","A145"), @("B128","B145"), @("D128","D145"), @("E128","E145"), @("F128","F145"),
  @("A129","A146"), @("B129","B146"), @("D129","D146"), @("E129","E146"), @("F129","F146"),
  @("A130","A147"), @("B130","B147"),
  @("A131","A148"), @("B131","B148"), @("D131","D148"), @("E131","E148"), @("F131","F148"),
  @("A133","A150"),
  @("A134","A151"), @("B134","B151"),
  @("A135","A152"), @("B135","B152"), @("D135","D152"), @("E135","E152"), @("F135","F152"),
  @("A136","A153"), @("B136","B153"), @("D136","D153"), @("E136","E153"), @("F136","F153"),
  @("F137","F154")
)

foreach ($pair in $cellMap) {
  $ws.Range($pair[0]).Copy()
  $ws.Range($pair[1]).PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = $false

# Row 139: commit label
$ws.Range("A139").Value = "Commit 8"

# Row 140: section title ("MARS Tool Output" / "Calulations")
$ws.Range("A140").Value = "MARS Tool Output"
$ws.Range("D140").Value = "Calulations"

# Row 142: "Instruction Statistics Tool"
$ws.Range("A142").Value = "Instruction Statistics Tool"

# Row 143: table header
$ws.Range("A143").Value = "Instruction type"
$ws.Range("B143").Value = "Count"
$ws.Range("D143").Value = "Adjusted count"
$ws.Range("E143").Value = "CPI"
$ws.Range("F143").Value = "Total cycles"

# Row 144: ALU
$ws.Range("A144").Value = "ALU"
$ws.Range("B144").Value = 2943
$ws.Range("D144").Formula = "=B144"
$ws.Range("E144").Value = 1
$ws.Range("F144").Formula = "=D144*E144"

# Row 145: Jump
$ws.Range("A145").Value = "Jump"
$ws.Range("B145").Value = 72
$ws.Range("D145").Formula = "=B145"
$ws.Range("E145").Value = 1
$ws.Range("F145").Formula = "=D145*E145"

# Row 146: Branch
$ws.Range("A146").Value = "Branch"
$ws.Range("B146").Value = 915
$ws.Range("D146").Formula = "=B146"
$ws.Range("E146").Value = 2
$ws.Range("F146").Formula = "=D146*E146"

# Row 147: Memory
$ws.Range("A147").Value = "Memory"
$ws.Range("B147").Value = 617

# Row 148: Other
$ws.Range("A148").Value = "Other"
$ws.Range("B148").Value = 712
$ws.Range("D148").Formula = "=B148-(B152+B153-B147)"
$ws.Range("E148").Value = 5
$ws.Range("F148").Formula = "=D148*E148"

# Row 150: "Data Cache Simulation Tool"
$ws.Range("A150").Value = "Data Cache Simulation Tool"

# Row 151: table header
$ws.Range("A151").Value = "Access"
$ws.Range("B151").Value = "Count"

# Row 152: Cache hit
$ws.Range("A152").Value = "Cache hit"
$ws.Range("B152").Value = 599
$ws.Range("D152").Formula = "=B152"
$ws.Range("E152").Value = 2
$ws.Range("F152").Formula = "=D152*E152"

# Row 153: Cache miss
$ws.Range("A153").Value = "Cache miss"
$ws.Range("B153").Value = 110
$ws.Range("D153").Formula = "=B153"
$ws.Range("E153").Value = 40
$ws.Range("F153").Formula = "=D153*E153"

# Row 154: total
$ws.Range("F154").Formula = "=SUM(F144:F153)"

# --- Update the sheet view to reflect the newly-added rows.
$ws.Application.ActiveWindow.ScrollRow = 137
$ws.Range("C154").Select()
